$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.14763442785013
$ws.Range("C2").Value = 5.674338765439493
$ws.Range("D2").Value = 6.995812154052948
$ws.Range("E2").Value = 10.89328524447621
$ws.Range("F2").Value = 36.54564082404119
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 12.54395826449728
$ws.Range("M2").Value = 15.45836519353912
$ws.Range("N2").Value = 22.05602951071957
$ws.Range("B3").Value = 12.90300152081777
$ws.Range("C3").Value = 5.468126719357534
$ws.Range("D3").Value = 7.004714139653245
$ws.Range("E3").Value = 10.68004866482394
$ws.Range("F3").Value = 36.32024241319986
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 12.38199045569924
$ws.Range("M3").Value = 15.3057058912435
$ws.Range("N3").Value = 22.0791869710074
$ws.Range("B4").Value = 12.75462639527421
$ws.Range("C4").Value = 5.339488557428469
$ws.Range("D4").Value = 7.010280847802146
$ws.Range("E4").Value = 10.55058114629494
$ws.Range("F4").Value = 36.18938287075239
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 12.28525165783836
$ws.Range("M4").Value = 15.21569009320429
$ws.Range("N4").Value = 22.09520667110844
$ws.Range("B5").Value = 12.69471338035311
$ws.Range("C5").Value = 5.286654779910362
$ws.Range("D5").Value = 7.012574888963574
$ws.Range("E5").Value = 10.49826652806886
$ws.Range("F5").Value = 36.13798917063363
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 12.24656038395257
$ws.Range("M5").Value = 15.17998091453097
$ws.Range("N5").Value = 22.10218718289555
$ws.Range("B6").Value = 12.68480080359064
$ws.Range("C6").Value = 5.277859812751306
$ws.Range("D6").Value = 7.012957362539182
$ws.Range("E6").Value = 10.48960882402887
$ws.Range("F6").Value = 36.12957304756824
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 12.24018125590398
$ws.Range("M6").Value = 15.17411127187021
$ws.Range("N6").Value = 22.10337359759234
$ws.Range("B7").Value = 12.75381603632765
$ws.Range("C7").Value = 5.338777560639823
$ws.Range("D7").Value = 7.010311682299666
$ws.Range("E7").Value = 10.54987371151472
$ws.Range("F7").Value = 36.18868188506229
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 12.28472683205448
$ws.Range("M7").Value = 15.21520451832634
$ws.Range("N7").Value = 22.09529898191569
$ws.Range("B8").Value = 13.06295948568585
$ws.Range("C8").Value = 5.603711838228624
$ws.Range("D8").Value = 6.998860750246856
$ws.Range("E8").Value = 10.81950435178496
$ws.Range("F8").Value = 36.46638290454301
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 12.48757919285895
$ws.Range("M8").Value = 15.40498020206297
$ws.Range("N8").Value = 22.06364012056163
$ws.Range("B9").Value = 13.67975060619454
$ws.Range("C9").Value = 6.103395136797101
$ws.Range("D9").Value = 6.977195878263608
$ws.Range("E9").Value = 11.35647424617703
$ws.Range("F9").Value = 37.06894966650856
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 12.9046529415831
$ws.Range("M9").Value = 15.80485458249603
$ws.Range("N9").Value = 22.01586984991549
$ws.Range("B10").Value = 14.13414272667927
$ws.Range("C10").Value = 6.453910265604036
$ws.Range("D10").Value = 6.961745832631292
$ws.Range("E10").Value = 11.751591151212
$ws.Range("F10").Value = 37.54440972837766
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 13.21979691446712
$ws.Range("M10").Value = 16.1130572421851
$ws.Range("N10").Value = 21.98953096584893
$ws.Range("B11").Value = 14.34007259902631
$ws.Range("C11").Value = 6.608974829726807
$ws.Range("D11").Value = 6.954815415645907
$ws.Range("E11").Value = 11.93058177620847
$ws.Range("F11").Value = 37.76721417899339
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 13.36440868978306
$ws.Range("M11").Value = 16.25584421256912
$ws.Range("N11").Value = 21.97945724463399
$ws.Range("B12").Value = 14.41785822520792
$ws.Range("C12").Value = 6.667007435208556
$ws.Range("D12").Value = 6.95220488751976
$ws.Range("E12").Value = 11.99818295174618
$ws.Range("F12").Value = 37.85246691022461
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 13.41929567371657
$ws.Range("M12").Value = 16.31023829092713
$ws.Range("N12").Value = 21.9759174692872
$ws.Range("B13").Value = 14.40111569064729
$ws.Range("C13").Value = 6.654540450158086
$ws.Range("D13").Value = 6.952766497845496
$ws.Range("E13").Value = 11.98363286583399
$ws.Range("F13").Value = 37.83406790290318
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 13.40747004538671
$ws.Range("M13").Value = 16.29850991687976
$ws.Range("N13").Value = 21.97666758715949
$ws.Range("B14").Value = 14.3464764552974
$ws.Range("C14").Value = 6.613763272442482
$ws.Range("D14").Value = 6.954600369124641
$ws.Range("E14").Value = 11.93614733975243
$ws.Range("F14").Value = 37.77421062872055
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 13.36892206433735
$ws.Range("M14").Value = 16.26031303211669
$ws.Range("N14").Value = 21.97916051215164
$ws.Range("B15").Value = 14.31298045481066
$ws.Range("C15").Value = 6.588695059514549
$ws.Range("D15").Value = 6.955725468449174
$ws.Range("E15").Value = 11.90703577216964
$ws.Range("F15").Value = 37.73765944736961
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 13.34532508085136
$ws.Range("M15").Value = 16.23695707720418
$ws.Range("N15").Value = 21.98072332026368
$ws.Range("B16").Value = 14.12066184636213
$ws.Range("C16").Value = 6.443683198650884
$ws.Range("D16").Value = 6.962200703847016
$ws.Range("E16").Value = 11.73987246422129
$ws.Range("F16").Value = 37.52997534500992
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 13.21036671152932
$ws.Range("M16").Value = 16.10377369116572
$ws.Range("N16").Value = 21.99022774876619
$ws.Range("B17").Value = 14.0024209568775
$ws.Range("C17").Value = 6.353557385950142
$ws.Range("D17").Value = 6.966197976021848
$ws.Range("E17").Value = 11.63707957538842
$ws.Range("F17").Value = 37.40419777405428
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 13.12785569017474
$ws.Range("M17").Value = 16.02269787524393
$ws.Range("N17").Value = 21.99654748895108
$ws.Range("B18").Value = 13.93434372916196
$ws.Range("C18").Value = 6.301309415409522
$ws.Range("D18").Value = 6.96850632711534
$ws.Range("E18").Value = 11.57788925174082
$ws.Range("F18").Value = 37.33247052824379
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 13.08051850428406
$ws.Range("M18").Value = 15.97631131828421
$ws.Range("N18").Value = 22.0003619911176
$ws.Range("B19").Value = 13.91128486317209
$ws.Range("C19").Value = 6.28355056287668
$ws.Range("D19").Value = 6.969289485923931
$ws.Range("E19").Value = 11.55783920730942
$ws.Range("F19").Value = 37.30829248938076
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 13.06451339186419
$ws.Range("M19").Value = 15.96064935015932
$ws.Range("N19").Value = 22.0016843372971
$ws.Range("B20").Value = 14.0150155810573
$ws.Range("C20").Value = 6.363194293991817
$ws.Range("D20").Value = 6.965771506257953
$ws.Range("E20").Value = 11.64802946828529
$ws.Range("F20").Value = 37.41752358228732
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 13.1366270053792
$ws.Range("M20").Value = 16.03130340894742
$ws.Range("N20").Value = 21.99585615459454
$ws.Range("B21").Value = 14.36253127990622
$ws.Range("C21").Value = 6.625759575905206
$ws.Range("D21").Value = 6.9540613418786
$ws.Range("E21").Value = 11.9501003860967
$ws.Range("F21").Value = 37.79176868658914
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 13.38024155859632
$ws.Range("M21").Value = 16.27152397840823
$ws.Range("N21").Value = 21.97842081390433
$ws.Range("B22").Value = 14.58847703930421
$ws.Range("C22").Value = 6.793335638958403
$ws.Range("D22").Value = 6.946488823734172
$ws.Range("E22").Value = 12.14644939430503
$ws.Range("F22").Value = 38.04147085826292
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 13.54016567869123
$ws.Range("M22").Value = 16.43038686494736
$ws.Range("N22").Value = 21.96862852710774
$ws.Range("B23").Value = 14.46801982303436
$ws.Range("C23").Value = 6.704282274616173
$ws.Range("D23").Value = 6.950523098551735
$ws.Range("E23").Value = 12.04177476212747
$ws.Range("F23").Value = 37.90775113728164
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 13.45476405442365
$ws.Range("M23").Value = 16.34544393936217
$ws.Range("N23").Value = 21.97370802562524
$ws.Range("B24").Value = 14.00932185327432
$ws.Range("C24").Value = 6.358838793610628
$ws.Range("D24").Value = 6.965964281225675
$ws.Range("E24").Value = 11.64307930699894
$ws.Range("F24").Value = 37.41149716117736
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 13.13266117970099
$ws.Range("M24").Value = 16.02741214265785
$ws.Range("N24").Value = 21.99616814243779
$ws.Range("B25").Value = 13.51230730683005
$ws.Range("C25").Value = 5.970839404178887
$ws.Range("D25").Value = 6.982973701827112
$ws.Range("E25").Value = 11.21079204752995
$ws.Range("F25").Value = 36.89998846239462
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 12.79007103736927
$ws.Range("M25").Value = 15.69396148765037
$ws.Range("N25").Value = 22.02725728114239
